$d = $word.ActiveDocument

# Locate the two list-item paragraphs that contain the screenshots
# ("Fix the output of N best moves." and "Fix the error getting
# evaluation.") and remove them (along with their inline images),
# leaving the trailing empty list paragraph intact.
$target1 = $null
$target2 = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    if ($text -like "Fix the output of N best moves.*") {
        $target1 = $p
    }
    elseif ($text -like "Fix the error getting evaluation.*") {
        $target2 = $p
    }
}

if ($target1 -ne $null -and $target2 -ne $null) {
    $rng = $d.Range($target1.Range.Start, $target2.Range.End)
    $rng.Delete()
}
